$d = $word.ActiveDocument

# 1) Table-cell paragraphs: the generic Paragraphs(...).LineSpacingRule
#    setter does not reach paragraphs nested inside table cells in this
#    runtime, so rebuild each cell paragraph's XML with the spacing
#    element added directly.
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables($ti)
    $rowCount = $t.Rows.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        $colCount = $t.Rows($r).Cells.Count
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $t.Cell($r, $c)
            $para = $cell.Range.Paragraphs(1)
            $text = $para.Range.Text
            $text = $text.TrimEnd([char]13, [char]7)
            $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
            $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:line='240' w:lineRule='auto'/></w:pPr><w:r><w:t>$escaped</w:t></w:r></w:p>"
            $null = $para.Range.InsertXML($xml)
        }
    }
}

# 2) All remaining (non-table) body paragraphs: set single line spacing
#    (w:spacing w:line="240" w:lineRule="auto") via the standard object
#    model property.
$d.Paragraphs.LineSpacingRule = 0
